# Applies the "added customization to the rev bar graphics" change:
#  - rows 8-17 col E: value becomes the generic placeholder "<value>"
#  - rows 47-56: new rows with SingleUseId60..69 / small / Left / 1..10 / LTR

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# 1) Rows 8-17, column E -> "<value>"
for ($r = 8; $r -le 17; $r++) {
    $ws.Range("E$r").Value = "<value>"
}

# 2) New rows 47-56
$ids = @(60,61,62,63,64,65,66,67,68,69)
$nums = @(1,2,3,4,5,6,7,8,9,10)

for ($i = 0; $i -lt 10; $i++) {
    $row = 47 + $i
    $ws.Range("B$row").Value = "SingleUseId$($ids[$i])"
    $ws.Range("C$row").Value = "small"
    $ws.Range("D$row").Value = "Left"
    # The E column holds the TextId number as text (not a numeric value),
    # so force text storage before writing, matching the rest of the sheet.
    $eCell = $ws.Range("E$row")
    $eCell.NumberFormat = "@"
    $eCell.Value = "$($nums[$i])"
    $ws.Range("F$row").Value = "LTR"
}
